$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.250.27'
$ws.Range("E2").Value = '  -1.10%  '
$ws.Range("D3").Value = '3.069.38'
$ws.Range("E3").Value = '  -1.73%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '574.46'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.99%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '169.69'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.65%  '
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("D8").Value = '3.066.21'
$ws.Range("E8").Value = '  -1.70%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.510'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.20%  '
$ws.Range("E10").Value = '  -2.04%  '
$ws.Range("E11").Value = '  -2.70%  '
$ws.Range("E12").Value = '  -3.12%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000238'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.81%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.57'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -4.29%  '
$ws.Range("E15").Value = '  -1.50%  '
$ws.Range("D16").Value = '3.581.31'
$ws.Range("E16").Value = '  -1.63%  '
$ws.Range("D17").Value = '66.210.39'
$ws.Range("E17").Value = '  -1.11%  '
$ws.Range("D19").Value = '3.070.49'
$ws.Range("E19").Value = '  -1.70%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '16.56'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.62%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '484.27'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.74%  '
$ws.Range("B22").Value = 'Uniswap'
$ws.Range("C22").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.64'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.80%  '
$ws.Range("B23").Value = 'Polygon'
$ws.Range("C23").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.683'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.78%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '82.16'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.04%  '
$ws.Range("E25").Value = '  -4.71%  '
$ws.Range("E26").Value = '  -3.61%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.04'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.77%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.999'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.16%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.78'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.81%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.23'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -5.50%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.59'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.70%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '27.58'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.60%  '
$ws.Range("E33").Value = '  -3.22%  '
$ws.Range("E34").Value = '  -3.78%  '
$ws.Range("E35").Value = '  -0.02%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '48.06'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.34%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.941'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.32%  '
$ws.Range("E38").Value = '  -5.11%  '
$ws.Range("E39").Value = '  -1.68%  '
$ws.Range("E40").Value = '  -3.84%  '
$ws.Range("E41").Value = '  -5.22%  '
$ws.Range("D43").Value = '2.773.39'
$ws.Range("E43").Value = '  -1.78%  '
$ws.Range("E44").Value = '  -1.00%  '
$ws.Range("E45").Value = '  -2.89%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '134.42'
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '364.69'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.62%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '24.11'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.46%  '
$ws.Range("E50").Value = '  -2.70%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.106'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.36%  '
